# Update cryptocurrency price/volume snapshot (scraper refresh).
# D column = Price, E column = Volume(1h) % change, both stored as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.177.09"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.319.20"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'309.93"
$ws.Range("E5").Value = "  +2.82%  "
$ws.Range("D6").Value = "'108.74"
$ws.Range("E6").Value = "  -4.61%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "'43.91"
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "'8.82"
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("E13").Value = "  +17.78%  "
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "2.661.49"
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").Value = "2.317.57"
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").Value = "43.158.04"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "'7.27"
$ws.Range("E20").Value = "  -5.98%  "
$ws.Range("D21").Value = "'76.00"
$ws.Range("E21").Value = "  +3.31%  "
$ws.Range("E22").Value = "  -5.68%  "
$ws.Range("D23").Value = "'2.53"
$ws.Range("E23").Value = "  +7.69%  "
$ws.Range("D24").Value = "'254.69"
$ws.Range("E24").Value = "  +9.35%  "
$ws.Range("D25").Value = "'9.05"
$ws.Range("E25").Value = "  -4.50%  "
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'39.20"
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").Value = "'22.52"
$ws.Range("E30").Value = "  +5.67%  "
$ws.Range("D31").Value = "'173.80"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").Value = "'3.16"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").Value = "'0.0906"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'5.78"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("D35").Value = "'5.03"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("E37").Value = "  -6.89%  "
$ws.Range("D38").Value = "'0.0377"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "'2.72"
$ws.Range("E40").Value = "  +4.69%  "
$ws.Range("D41").Value = "'1.49"
$ws.Range("E41").Value = "  +11.09%  "
$ws.Range("D42").Value = "'71.76"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'12.48"
$ws.Range("E45").Value = "  -6.83%  "
$ws.Range("D46").Value = "'5.73"
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("D47").Value = "'110.74"
$ws.Range("E47").Value = "  +4.39%  "
$ws.Range("D48").Value = "'9.11"
$ws.Range("E48").Value = "  +4.25%  "
$ws.Range("E49").Value = "  -5.94%  "
$ws.Range("D50").Value = "'0.0991"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "'71.16"
$ws.Range("E51").Value = "  +2.76%  "
